$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 364; Excel shifts rows 364..441 down to 365..442
# and picks up formatting (e.g. the date-formatted style on column D) from the
# row insert, matching the target workbook's <dimension ref="A1:R442"/>.
$ws.Rows.Item(364).Insert()

# Populate the newly inserted row 364 with the new data record.
$ws.Range("A364").Value = 9
$ws.Range("B364").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C364").Value = "Metropolitana"
$ws.Range("D364").Value = 44722
$ws.Range("E364").Value = 13
$ws.Range("F364").Value = 100112012
$ws.Range("G364").Value = "Espinaca"
$ws.Range("H364").Value = "Sin especificar"
$ws.Range("I364").Value = "Primera"
$ws.Range("J364").Value = 160
$ws.Range("K364").Value = 6000
$ws.Range("L364").Value = 7000
$ws.Range("M364").Value = 6500
$ws.Range("N364").Value = "`$/cuna 10 kilos"
$ws.Range("O364").Value = "Provincia de Chacabuco"
$ws.Range("P364").Value = 650
$ws.Range("Q364").Value = 10
$ws.Range("R364").Value = "Hortaliza"
